$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a blank row above row 18; everything from 18 downward shifts to 19+.
$ws.Rows.Item(18).Insert(-4121)   # xlShiftDown

# Row 17 carries exactly the per-column formatting (borders/fonts/number
# format) the new row 18 needs, so stamp it onto the freshly inserted row.
$ws.Range("A17:H17").Copy()
$ws.Range("A18:H18").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows.Item(18).RowHeight = 16.5

# New task row: Yale's 1.1/2.0 CAS plugin work.
$ws.Range("A18").Value = "1.1/2.0"
$ws.Range("B18").Value = "Yale"
$ws.Range("D18").Value = 5368
$ws.Range("E18").Value = "low"
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = "Implement a CAS plugin for the Authentication Service."

# Match the author's cursor position after the edit.
$ws.Range("J8").Select()
